$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Engineering"
$ws.Range("A2").Value = "Have you ever been to a work shop?"

$ws.Range("A3").Value = "Have you used auto card before?"
$ws.Range("B3").Value = "Engineering"

$ws.Range("A4").Value = "Have u fixed a tire before"
$ws.Range("B4").Value = "Engineering"

$ws.Range("C4").Select()
